$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 151
$ws1.Range("F5").Value = 3332
$ws1.Range("F6").Value = 1084
$ws1.Range("F7").Value = 2219
$ws1.Range("F8").Value = 2130
$ws1.Range("F9").Value = 1115
$ws1.Range("F16").Value = 98
$ws1.Range("F17").Value = 219
$ws1.Range("F18").Value = 1590
$ws1.Range("F19").Value = 640
$ws1.Range("F20").Value = 729
$ws1.Range("F21").Value = 612
$ws1.Range("F22").Value = 12280
$ws1.Range("F23").Value = 12335
$ws1.Range("F25").Value = 705
$ws1.Range("F27").Value = 40
$ws1.Range("F29").Value = 376
$ws1.Range("F30").Value = 1924
$ws1.Range("F34").Value = 594

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 1

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 151
$ws4.Range("F6").Value = 3332
$ws4.Range("F7").Value = 1084
$ws4.Range("F8").Value = 2219
$ws4.Range("F9").Value = 2130
$ws4.Range("F10").Value = 1115
$ws4.Range("F15").Value = 1
$ws4.Range("F19").Value = 98
$ws4.Range("F21").Value = 219
$ws4.Range("F22").Value = 1590
$ws4.Range("F23").Value = 640
$ws4.Range("F24").Value = 729
$ws4.Range("F25").Value = 612
$ws4.Range("F26").Value = 12280
$ws4.Range("F27").Value = 12335
$ws4.Range("F29").Value = 705
$ws4.Range("F31").Value = 40
$ws4.Range("F33").Value = 376
$ws4.Range("F34").Value = 1924
$ws4.Range("F40").Value = 594
